$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Max-Min HR Window 10 Sec Stats")

# Define name used by the new formulas
$wb.Names.Add("CalculationArea", '=''Max-Min HR Window 10 Sec Stats''!$B$2:$B$26')

# Row1 new header + new static value
$ws3.Range("G1").Value = "3x Median"
$ws3.Range("M1").Value = 2

# Row2
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Formula = "=IF(A2<30, A2, 30)"
$ws3.Range("C2").Formula = "=MEDIAN(CalculationArea)"
$ws3.Range("D2").Formula = '=COUNTIF(CalculationArea, "<="&C2)/COUNT(CalculationArea)'
$ws3.Range("E2").Formula = '=COUNTIF(CalculationArea, ">"&C2)/COUNT(CalculationArea)-F2'
$ws3.Range("F2").Formula = '=COUNTIF(CalculationArea, ">="&G2)/COUNT(CalculationArea)'
$ws3.Range("G2").Formula = "=3*C2"
$ws3.Range("M2").Value = 4

# Row3
$ws3.Range("A3").Value = 4
$ws3.Range("B3").Formula = "=IF(A3<30, A3, 30)"
$ws3.Range("M3").Value = 4

# Row4
$ws3.Range("A4").Value = 4
$ws3.Range("B4").Formula = "=IF(A4<30, A4, 30)"
$ws3.Range("M4").Value = 5

# Row5
$ws3.Range("A5").Value = 10
$ws3.Range("B5").Formula = "=IF(A5<30, A5, 30)"
$ws3.Range("M5").Value = 5

# Row6
$ws3.Range("A6").Value = 14
$ws3.Range("B6").ClearContents()
$ws3.Range("M6").Value = 7

# Row7
$ws3.Range("A7").Value = 12
$ws3.Range("B7").ClearContents()
$ws3.Range("M7").Value = 7

# Row8
$ws3.Range("A8").Value = 6
$ws3.Range("B8").ClearContents()
$ws3.Range("M8").Value = 12

# Row9
$ws3.Range("A9").Value = 7
$ws3.Range("B9").ClearContents()

# Row10
$ws3.Range("A10").Value = 3
$ws3.Range("B10").ClearContents()

# Row11
$ws3.Range("A11").Value = 3
$ws3.Range("B11").ClearContents()

# Row12
$ws3.Range("A12").Value = 9
$ws3.Range("B12").ClearContents()

# Row13
$ws3.Range("A13").Value = 13
$ws3.Range("B13").ClearContents()

# Row14
$ws3.Range("A14").Value = 22
$ws3.Range("B14").ClearContents()

# Row15
$ws3.Range("A15").Value = 6
$ws3.Range("B15").ClearContents()

# Row16
$ws3.Range("A16").Value = 7
$ws3.Range("B16").ClearContents()

# Row17
$ws3.Range("A17").Value = 3
$ws3.Range("B17").ClearContents()

# Row18
$ws3.Range("A18").Value = 7
$ws3.Range("B18").ClearContents()

# Row19
$ws3.Range("A19").Value = 3
$ws3.Range("B19").ClearContents()

# Row20
$ws3.Range("A20").Value = 4
$ws3.Range("B20").ClearContents()

# Row21
$ws3.Range("A21").Value = 17
$ws3.Range("B21").ClearContents()

# Row22
$ws3.Range("A22").Value = 20
$ws3.Range("B22").ClearContents()

# Row23
$ws3.Range("A23").Value = 25
$ws3.Range("B23").ClearContents()

# Row24
$ws3.Range("A24").Value = 87
$ws3.Range("B24").ClearContents()

# Row25
$ws3.Range("A25").Value = 4
$ws3.Range("B25").ClearContents()

# Row26
$ws3.Range("A26").Value = 44
$ws3.Range("B26").ClearContents()

Write-Output "done"
